# Apply an automatic data update to rows 10-13: the observation records
# (Id, Lokalnamn, Ost, Nord and the public comment) are cyclically
# reassigned among the four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for the columns that change: A (Id), P (Lokalnamn),
# Q (Ost), R (Nord), AC (Publik kommentar).
$rows = 10..13
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        A  = $ws.Range("A$r").Value2
        P  = $ws.Range("P$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        AC = $ws.Range("AC$r").Value2
    }
}

# Each row's new data comes from the next row (row 13 wraps to row 10).
foreach ($r in $rows) {
    $srcRow = if ($r -eq 13) { 10 } else { $r + 1 }
    $src = $data[$srcRow]

    $ws.Range("A$r").Value = $src.A
    $ws.Range("P$r").Value = $src.P
    $ws.Range("Q$r").Value = $src.Q
    $ws.Range("R$r").Value = $src.R

    if ($src.AC) {
        $ws.Range("AC$r").Value = $src.AC
    } else {
        $ws.Range("AC$r").Value = ""
    }
}

$wb.Save()
